$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new L2 cell with value "citrsi" (previously absent, now inserted before M2)
$ws.Range("L2").Value = "citrsi"

# Update tnrsdate column (T2:T11) from 45909 to the new run date 45979
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 20).Value = 45979
}
